# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / Handoff / Handback datetime
# stamps on each sheet to reflect the new report-generation timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-28 21:03:26"

$zhcn.Range("H2").Value = "2016-08-28 21:03:20"
$zhcn.Range("K2").Value = "2016-08-28 21:03:48"

$dede.Range("H2").Value = "2016-08-28 21:03:26"
$dede.Range("K2").Value = "2016-08-28 21:03:55"
